$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 corresponds to the "serialVersion" property.
# ControlledBy (column C) changes from MN to CN, and
# ModifiableBy (column F) changes from "MN-service-subject" to
# "CN replication processes" since serialVersion is now managed by the
# CN replication processes.
$ws.Range("C2").Value = "CN"
$ws.Range("F2").Value = "CN replication processes"

# Update the active selection on the sheet to E7.
$ws.Range("E7").Select()

$wb.Save()
